$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 2, columns A-R (1-18) with 0, mirroring the header row's extent
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(2, $col).Value = 0
}

# Move the active selection to A3, as if the user tabbed/entered through row 2
$ws.Range("A3").Select()
